{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the content of `async (context) => { ... }`.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that introduces \"Problem 3\" (Predicting Fingers) so\n// we anchor relative to content rather than a fixed index.\nlet problem3Index = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Predicting Fingers\") !== -1) {\n    problem3Index = i;\n    break;\n  }\n}\nif (problem3Index === -1) {\n  throw new Error(\"Could not locate the 'Predicting Fingers' (Problem 3) heading paragraph.\");\n}\n\n// Find the bookmark-only (\"_GoBack\") paragraph that follows Problem 3's\n// description text \u2014 the new paragraph must land immediately before it.\nlet goBackIndex = -1;\nfor (let i = problem3Index + 1; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Define the problem\") !== -1) {\n    break;\n  }\n  if (paragraphs.items[i].text === \"\" && i > problem3Index + 1) {\n    goBackIndex = i;\n  }\n}\nif (goBackIndex === -1) {\n  throw new Error(\"Could not locate the target blank paragraph before the numbered list.\");\n}\n\n// Insert a new (empty) paragraph directly before that blank paragraph; its\n// content will hold the constraints / sub-goal text for Problem 3.\nconst anchorParagraph = paragraphs.items[goBackIndex];\nconst newParagraph = anchorParagraph.insertParagraph(\"\", \"Before\");\n\n// The source splits the sentence across two runs (\"The co\" + \"nstraints\n// are...\"). Office.js's insertText() coalesces same-formatted adjacent runs,\n// so use insertOoxml (Flat OPC) on the new paragraph's content range to get\n// the exact two-run structure.\nconst firstRun = \"The co\";\nconst secondRun = \"nstraints are that the little girl has used a system that has difficulty testing and seeing from a simple demonstration without the potential of an error occurring.   The sub-goal will be to come up with a method of applying her system to an equation with variables to simplify toward the solution.\";\n\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>' + firstRun + '</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">' + secondRun + '</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nconst contentRange = newParagraph.getRange(\"Content\");\ncontentRange.insertOoxml(flatOpc, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that introduces \"Problem 3\" (Predicting Fingers) so\n# we anchor relative to content rather than a fixed, brittle index.\n$problem3Index = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Predicting Fingers*\") {\n        $problem3Index = $i\n        break\n    }\n}\nif ($problem3Index -eq -1) {\n    throw \"Could not locate the 'Predicting Fingers' (Problem 3) heading paragraph.\"\n}\n\n# Starting from there, find the blank paragraph that sits immediately before\n# the \"_GoBack\"-bookmark-only paragraph (i.e., the blank line right above\n# \"1) Define the problem\"). That is where the new paragraph must land.\n$goBackIndex = -1\nfor ($i = $problem3Index + 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Define the problem*\") {\n        break\n    }\n    if ($t -eq [char]13) {\n        $goBackIndex = $i\n    }\n}\nif ($goBackIndex -eq -1) {\n    throw \"Could not locate the target blank paragraph before the numbered list.\"\n}\n\n# Insert a new (empty) paragraph right before that blank paragraph.\n$anchorRange = $d.Paragraphs.Item($goBackIndex - 1).Range\n$anchorRange.Collapse(0) # wdCollapseEnd\n$anchorRange.InsertParagraphAfter() | Out-Null\n\n# The new paragraph is now at $goBackIndex. The source splits the sentence\n# across two runs (\"The co\" + \"nstraints are...\"). Plain InsertAfter() calls\n# get coalesced into a single run on save, so use InsertXML (Flat OPC) to\n# seed the paragraph with the exact two-run structure.\n$newParaRange = $d.Paragraphs.Item($goBackIndex).Range\n\n$firstRun = \"The co\"\n$secondRun = \"nstraints are that the little girl has used a system that has difficulty testing and seeing from a simple demonstration without the potential of an error occurring.   The sub-goal will be to come up with a method of applying her system to an equation with variables to simplify toward the solution.\"\n\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:t>' + $firstRun + '</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">' + $secondRun + '</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$newParaRange.InsertXML($flatOpc) | Out-Null\n"}
